$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("B1").Value = "GDP_Base"

# Delete the now-redundant "High_Renewables" column (old column C, duplicate of B)
$ws.Range("C:C").Delete()

# Rename the remaining (shifted) header, old D1 -> now C1
$ws.Range("C1").Value = "GDP_High_EconGrowth"

# Update selection to match saved workbook state
$ws.Range("C2").Select()
